$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 105.166664
$ws.Range("I9").Value = 95.25
$ws.Range("K9").Value = 95.25
$ws.Range("M9").Value = 73.75

$ws.Range("H17").Value = 952.4167
$ws.Range("J17").Value = 876.8
$ws.Range("L17").Value = 2630.4
$ws.Range("N17").Value = -2966.4

$ws.Range("H40").Value = 5258.3335
$ws.Range("I40").Value = 6973.4
$ws.Range("J40").Value = 3114.5
$ws.Range("K40").Value = 6973.4
$ws.Range("L40").Value = 3114.5
$ws.Range("M40").Value = -6798.4
$ws.Range("N40").Value = -3464.5

$ws.Range("H64").Value = 5399.6
$ws.Range("J64").Value = 5666.3335
$ws.Range("L64").Value = 5666.3335
$ws.Range("N64").Value = -6162.3335

$ws.Range("H67").Value = 5399.6
$ws.Range("J67").Value = 5666.3335
$ws.Range("L67").Value = 5666.3335
$ws.Range("N67").Value = -7382.3335

$ws.Range("H92").Value = 1654.5
$ws.Range("I92").Value = 1259.4546
$ws.Range("K92").Value = 1259.4546
$ws.Range("M92").Value = -11.45460000000003

$ws.Range("H98").Value = 2530.6
$ws.Range("I98").Value = 2883.1875
$ws.Range("K98").Value = 2883.1875
$ws.Range("M98").Value = -1385.1875

$ws.Range("H100").Value = 7307.826
$ws.Range("J100").Value = 7666.7334
$ws.Range("L100").Value = 7666.7334
$ws.Range("N100").Value = -8748.733400000001

$ws.Range("H122").Value = 2530.6
$ws.Range("I122").Value = 2883.1875
$ws.Range("K122").Value = 8649.5625
$ws.Range("M122").Value = -6199.5625

$ws.Range("H132").Value = 19257.932
$ws.Range("I132").Value = 1833.8462
$ws.Range("K132").Value = 5501.5386
$ws.Range("M132").Value = -2971.5386

$ws.Range("H137").Value = 6464129.5
$ws.Range("I137").Value = 10016199
$ws.Range("J137").Value = 5821.636
$ws.Range("K137").Value = 30048597
$ws.Range("L137").Value = 17464.908
$ws.Range("M137").Value = -30046047
$ws.Range("N137").Value = -22564.908

$ws.Range("H138").Value = 4463.0386
$ws.Range("I138").Value = 2950
$ws.Range("J138").Value = 5572.6
$ws.Range("K138").Value = 8850
$ws.Range("L138").Value = 16717.8
$ws.Range("M138").Value = -3710
$ws.Range("N138").Value = -26997.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 20907.666
$ws.Range("I2").Value = 25402.715
$ws.Range("J2").Value = 5175
$ws.Range("K2").Value = 25402.715
$ws.Range("L2").Value = 5175
$ws.Range("M2").Value = -25289.715
$ws.Range("N2").Value = -5401

$ws.Range("H110").Value = 3321.2654
$ws.Range("I110").Value = 3276.5588
$ws.Range("K110").Value = 3276.5588
$ws.Range("M110").Value = -1231.5588

$ws.Range("H116").Value = 20907.666
$ws.Range("I116").Value = 25402.715
$ws.Range("J116").Value = 5175
$ws.Range("K116").Value = 25402.715
$ws.Range("L116").Value = 5175
$ws.Range("M116").Value = -23108.715
$ws.Range("N116").Value = -9763

$ws.Range("H122").Value = 4725.143
$ws.Range("I122").Value = 4307.5884
$ws.Range("K122").Value = 12922.7652
$ws.Range("M122").Value = -10472.7652

$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 20907.666
$ws.Range("I3").Value = 25402.715
$ws.Range("J3").Value = 5175
$ws.Range("K3").Value = 25402.715
$ws.Range("L3").Value = 5175
$ws.Range("M3").Value = -25288.715
$ws.Range("N3").Value = -5403

$ws.Range("H5").Value = 706.75
$ws.Range("I5").Value = 631.8
$ws.Range("J5").Value = 831.6667
$ws.Range("K5").Value = 631.8
$ws.Range("L5").Value = 831.6667
$ws.Range("M5").Value = -518.8
$ws.Range("N5").Value = -1057.6667

$ws.Range("H86").Value = 7741.85
$ws.Range("I86").Value = 6888.2354
$ws.Range("K86").Value = 6888.2354
$ws.Range("M86").Value = -5765.2354

$ws.Range("H89").Value = 7741.85
$ws.Range("I89").Value = 6888.2354
$ws.Range("K89").Value = 34441.177
$ws.Range("M89").Value = -28825.177

$ws.Range("H110").Value = 98199.8
$ws.Range("J110").Value = 98199.8
$ws.Range("L110").Value = 98199.8
$ws.Range("N110").Value = -106379.8

$ws.Range("H134").Value = 17201.5
$ws.Range("I134").Value = 3403.2856
$ws.Range("J134").Value = 36519
$ws.Range("K134").Value = 10209.8568
$ws.Range("L134").Value = 109557
$ws.Range("M134").Value = -7674.856800000001
$ws.Range("N134").Value = -114627

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7556.7646
$ws.Range("I31").Value = 8677.385
$ws.Range("J31").Value = 3914.75
$ws.Range("K31").Value = 8677.385
$ws.Range("L31").Value = 3914.75
$ws.Range("M31").Value = -8382.385
$ws.Range("N31").Value = -4504.75

$ws.Range("H34").Value = 7556.7646
$ws.Range("I34").Value = 8677.385
$ws.Range("J34").Value = 3914.75
$ws.Range("K34").Value = 8677.385
$ws.Range("L34").Value = 3914.75
$ws.Range("M34").Value = -8475.385
$ws.Range("N34").Value = -4318.75

$ws.Range("H105").Value = 4715.364
$ws.Range("I105").Value = 4874.3335
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 4874.3335
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -3127.3335
$ws.Range("N105").Value = -7494

$ws.Range("H107").Value = 1137.375
$ws.Range("J107").Value = 1705.25
$ws.Range("L107").Value = 1705.25
$ws.Range("N107").Value = -5545.25

$ws.Range("H122").Value = 333328.3
$ws.Range("I122").Value = 396530.7
$ws.Range("K122").Value = 1189592.1
$ws.Range("M122").Value = -1187142.1

$ws.Range("H132").Value = 7850.5557
$ws.Range("I132").Value = 8199.764999999999
$ws.Range("K132").Value = 24599.295
$ws.Range("M132").Value = -22069.295

$ws.Range("H134").Value = 3822.724
$ws.Range("I134").Value = 3353.35
$ws.Range("J134").Value = 4865.778
$ws.Range("K134").Value = 10060.05
$ws.Range("L134").Value = 14597.334
$ws.Range("M134").Value = -7525.049999999999
$ws.Range("N134").Value = -19667.334

$ws.Range("H138").Value = 79769.234
$ws.Range("J138").Value = 79769.234
$ws.Range("L138").Value = 79769.234
$ws.Range("N138").Value = -90049.234

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 34233548
$ws.Range("I4").Value = 35414000
$ws.Range("K4").Value = 106242000
$ws.Range("M4").Value = -106241888

$ws.Range("H132").Value = 897.38464
$ws.Range("I132").Value = 943
$ws.Range("K132").Value = 8487
$ws.Range("M132").Value = -5957

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 1651.2858
$ws.Range("I6").Value = 416.66666
$ws.Range("K6").Value = 416.66666
$ws.Range("M6").Value = -303.66666

$ws.Range("H16").Value = 1651.2858
$ws.Range("I16").Value = 416.66666
$ws.Range("K16").Value = 416.66666
$ws.Range("M16").Value = -166.66666

$ws.Range("H102").Value = 2328.7297
$ws.Range("I102").Value = 2373
$ws.Range("K102").Value = 2373
$ws.Range("M102").Value = -751

$ws.Range("H122").Value = 5232.6
$ws.Range("I122").Value = 4788.5835
$ws.Range("J122").Value = 5898.625
$ws.Range("K122").Value = 14365.7505
$ws.Range("L122").Value = 17695.875
$ws.Range("M122").Value = -11915.7505
$ws.Range("N122").Value = -22595.875

$ws.Range("H132").Value = 9367.210999999999
$ws.Range("I132").Value = 9714.286
$ws.Range("K132").Value = 29142.858
$ws.Range("M132").Value = -26612.858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2183.5
$ws.Range("I46").Value = 1192.5
$ws.Range("K46").Value = 1192.5
$ws.Range("M46").Value = -1004.5

$ws.Range("H93").Value = 1160.125
$ws.Range("I93").Value = 1192.5238
$ws.Range("K93").Value = 1192.5238
$ws.Range("M93").Value = 55.47620000000006

$ws.Range("H122").Value = 2594.1177
$ws.Range("I122").Value = 2250
$ws.Range("K122").Value = 6750
$ws.Range("M122").Value = -4300

$ws.Range("H132").Value = 6300.8184
$ws.Range("I132").Value = 6541
$ws.Range("J132").Value = 3899
$ws.Range("K132").Value = 19623
$ws.Range("L132").Value = 11697
$ws.Range("M132").Value = -17093
$ws.Range("N132").Value = -16757

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1249.8334
$ws.Range("J13").Value = 1219.8
$ws.Range("L13").Value = 1219.8
$ws.Range("N13").Value = -1499.8

$ws.Range("H96").Value = 59938.332
$ws.Range("I96").Value = 128499.25
$ws.Range("K96").Value = 128499.25
$ws.Range("M96").Value = -127126.25

$ws.Range("H103").Value = 46601.668
$ws.Range("J103").Value = 48676.875
$ws.Range("L103").Value = 48676.875
$ws.Range("N103").Value = -51020.875

$ws.Range("H122").Value = 2545.3157
$ws.Range("I122").Value = 2430.7334
$ws.Range("K122").Value = 7292.2002
$ws.Range("M122").Value = -4842.2002

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""
